$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("!!_Schema")
$ws1.Unprotect()
$ws1.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8' date='2020-03-09 12:59:52'"
$ws1.Range("A2").Value = "!!ObjTables type='Schema' description='Table/model and column/attribute definitions' date='2020-03-09 12:59:52' objTablesVersion='0.0.8'"
$ws1.Protect()

$ws2 = $wb.Worksheets.Item("!!Compound")
$ws2.Unprotect()
$ws2.Range("A1").Value = "!!ObjTables type='Data' id='Compound' description='Compound' name='Compound' date='2020-03-09 12:59:52' objTablesVersion='0.0.8'"
$ws2.Protect()

$ws3 = $wb.Worksheets.Item("!!Model")
$ws3.Unprotect()
$ws3.Range("A1").Value = "!!ObjTables type='Data' id='Model' description='Model' name='Model' date='2020-03-09 12:59:52' objTablesVersion='0.0.8'"
$ws3.Protect()

$ws4 = $wb.Worksheets.Item("!!Reaction")
$ws4.Unprotect()
$ws4.Range("A1").Value = "!!ObjTables type='Data' id='Reaction' description='Reaction' name='Reaction' date='2020-03-09 12:59:52' objTablesVersion='0.0.8'"
$ws4.Protect()
